$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.277420333333333
$ws.Range("H2").Value = 3.832261
$ws.Range("I2").Value = 0.01913942624337554
$ws.Range("J2").Value = 0.01913942624337554
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 150.1098686666667
$ws.Range("N2").Value = 450.329606
$ws.Range("O2").Value = 0.7276622610660995
$ws.Range("P2").Value = 0.7276622610660997
$ws.Range("Q2").Value = 191.7533984687962
$ws.Range("R2").Value = 1725.780586219166
$ws.Range("S2").Value = 0.01392703817576249
$ws.Range("T2").Value = 0.01392703817576249

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.277420333333333
$ws.Range("H3").Value = 3.832261
$ws.Range("I3").Value = 0.01913942624337554
$ws.Range("J3").Value = 0.01913942624337554
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 42.32476666666667
$ws.Range("N3").Value = 126.9743
$ws.Range("O3").Value = 0.2051706239258123
$ws.Range("P3").Value = 0.2051706239258124
$ws.Range("Q3").Value = 54.06651754358889
$ws.Range("R3").Value = 486.5986578923
$ws.Range("S3").Value = 0.003926848023935427
$ws.Range("T3").Value = 0.003926848023935428

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.277420333333333
$ws.Range("H4").Value = 3.832261
$ws.Range("I4").Value = 0.01913942624337554
$ws.Range("J4").Value = 0.01913942624337554
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.940565666666667
$ws.Range("N4").Value = 14.821697
$ws.Range("O4").Value = 0.02394954586187395
$ws.Range("P4").Value = 0.02394954586187395
$ws.Range("Q4").Value = 6.311179040768557
$ws.Range("R4").Value = 56.800611366917
$ws.Range("S4").Value = 0.0004583805665856764
$ws.Range("T4").Value = 0.0004583805665856764

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.277420333333333
$ws.Range("H5").Value = 3.832261
$ws.Range("I5").Value = 0.01913942624337554
$ws.Range("J5").Value = 0.01913942624337554
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.915377333333334
$ws.Range("N5").Value = 26.746132
$ws.Range("O5").Value = 0.04321756914621411
$ws.Range("P5").Value = 0.04321756914621412
$ws.Range("Q5").Value = 11.38868428493911
$ws.Range("R5").Value = 102.498158564452
$ws.Range("S5").Value = 0.0008271594770919476
$ws.Range("T5").Value = 0.0008271594770919477

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("H6").Value = 136.341755
$ws.Range("I6").Value = 0.6809303864519871
$ws.Range("J6").Value = 0.6809303864519872
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 150.1098686666667
$ws.Range("N6").Value = 450.329606
$ws.Range("O6").Value = 0.7276622610660995
$ws.Range("P6").Value = 0.7276622610660997
$ws.Range("Q6").Value = 6822.080978944281
$ws.Range("R6").Value = 61398.72881049853
$ws.Range("S6").Value = 0.4954873446342659
$ws.Range("T6").Value = 0.4954873446342661

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("H7").Value = 136.341755
$ws.Range("I7").Value = 0.6809303864519871
$ws.Range("J7").Value = 0.6809303864519872
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.32476666666667
$ws.Range("N7").Value = 126.9743
$ws.Range("O7").Value = 0.2051706239258123
$ws.Range("P7").Value = 0.2051706239258124
$ws.Range("Q7").Value = 1923.544322432944
$ws.Range("R7").Value = 17311.8989018965
$ws.Range("S7").Value = 0.1397069122383987
$ws.Range("T7").Value = 0.1397069122383988

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("H8").Value = 136.341755
$ws.Range("I8").Value = 0.6809303864519871
$ws.Range("J8").Value = 0.6809303864519872
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.940565666666667
$ws.Range("N8").Value = 14.821697
$ws.Range("O8").Value = 0.02394954586187395
$ws.Range("P8").Value = 0.02394954586187395
$ws.Range("Q8").Value = 224.5351312286928
$ws.Range("R8").Value = 2020.816181058235
$ws.Range("S8").Value = 0.01630797351907541
$ws.Range("T8").Value = 0.01630797351907542

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("H9").Value = 136.341755
$ws.Range("I9").Value = 0.6809303864519871
$ws.Range("J9").Value = 0.6809303864519872
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.915377333333334
$ws.Range("N9").Value = 26.746132
$ws.Range("O9").Value = 0.04321756914621411
$ws.Range("P9").Value = 0.04321756914621412
$ws.Range("Q9").Value = 405.1793973712956
$ws.Range("R9").Value = 3646.61457634166
$ws.Range("S9").Value = 0.02942815606024705
$ws.Range("T9").Value = 0.02942815606024706

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.4966396666666666
$ws.Range("H10").Value = 1.489919
$ws.Range("I10").Value = 0.007441088905245192
$ws.Range("J10").Value = 0.007441088905245193
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 150.1098686666667
$ws.Range("N10").Value = 450.329606
$ws.Range("O10").Value = 0.7276622610660995
$ws.Range("P10").Value = 0.7276622610660997
$ws.Range("Q10").Value = 74.55051513799044
$ws.Range("R10").Value = 670.954636241914
$ws.Range("S10").Value = 0.005414599577584583
$ws.Range("T10").Value = 0.005414599577584585

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4966396666666666
$ws.Range("H11").Value = 1.489919
$ws.Range("I11").Value = 0.007441088905245192
$ws.Range("J11").Value = 0.007441088905245193
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 42.32476666666667
$ws.Range("N11").Value = 126.9743
$ws.Range("O11").Value = 0.2051706239258123
$ws.Range("P11").Value = 0.2051706239258124
$ws.Range("Q11").Value = 21.02015800907778
$ws.Range("R11").Value = 189.1814220817
$ws.Range("S11").Value = 0.001526692853376596
$ws.Range("T11").Value = 0.001526692853376596

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4966396666666666
$ws.Range("H12").Value = 1.489919
$ws.Range("I12").Value = 0.007441088905245192
$ws.Range("J12").Value = 0.007441088905245193
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.940565666666667
$ws.Range("N12").Value = 14.821697
$ws.Range("O12").Value = 0.02394954586187395
$ws.Range("P12").Value = 0.02394954586187395
$ws.Range("Q12").Value = 2.453680885838111
$ws.Range("R12").Value = 22.083127972543
$ws.Range("S12").Value = 0.0001782106999984511
$ws.Range("T12").Value = 0.0001782106999984512

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4966396666666666
$ws.Range("H13").Value = 1.489919
$ws.Range("I13").Value = 0.007441088905245192
$ws.Range("J13").Value = 0.007441088905245193
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.915377333333334
$ws.Range("N13").Value = 26.746132
$ws.Range("O13").Value = 0.04321756914621411
$ws.Range("P13").Value = 0.04321756914621412
$ws.Range("Q13").Value = 4.427730027034222
$ws.Range("R13").Value = 39.849570243308
$ws.Range("S13").Value = 0.0003215857742855607
$ws.Range("T13").Value = 0.0003215857742855609

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.52156333333333
$ws.Range("H14").Value = 58.56469
$ws.Range("I14").Value = 0.2924890983993922
$ws.Range("J14").Value = 0.2924890983993922
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 150.1098686666667
$ws.Range("N14").Value = 450.329606
$ws.Range("O14").Value = 0.7276622610660995
$ws.Range("P14").Value = 0.7276622610660997
$ws.Range("Q14").Value = 2930.379308134682
$ws.Range("R14").Value = 26373.41377321214
$ws.Range("S14").Value = 0.2128332786784866
$ws.Range("T14").Value = 0.2128332786784867

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.52156333333333
$ws.Range("H15").Value = 58.56469
$ws.Range("I15").Value = 0.2924890983993922
$ws.Range("J15").Value = 0.2924890983993922
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 42.32476666666667
$ws.Range("N15").Value = 126.9743
$ws.Range("O15").Value = 0.2051706239258123
$ws.Range("P15").Value = 0.2051706239258124
$ws.Range("Q15").Value = 826.2456130518889
$ws.Range("R15").Value = 7436.210517467
$ws.Range("S15").Value = 0.06001017081010161
$ws.Range("T15").Value = 0.06001017081010163

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.52156333333333
$ws.Range("H16").Value = 58.56469
$ws.Range("I16").Value = 0.2924890983993922
$ws.Range("J16").Value = 0.2924890983993922
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.940565666666667
$ws.Range("N16").Value = 14.821697
$ws.Range("O16").Value = 0.02394954586187395
$ws.Range("P16").Value = 0.02394954586187395
$ws.Range("Q16").Value = 96.44756556432556
$ws.Range("R16").Value = 868.0280900789299
$ws.Range("S16").Value = 0.007004981076214405
$ws.Range("T16").Value = 0.007004981076214407

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.52156333333333
$ws.Range("H17").Value = 58.56469
$ws.Range("I17").Value = 0.2924890983993922
$ws.Range("J17").Value = 0.2924890983993922
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.915377333333334
$ws.Range("N17").Value = 26.746132
$ws.Range("O17").Value = 0.04321756914621411
$ws.Range("P17").Value = 0.04321756914621412
$ws.Range("Q17").Value = 174.0421032532311
$ws.Range("R17").Value = 1566.37892927908
$ws.Range("S17").Value = 0.01264066783458956
$ws.Range("T17").Value = 0.01264066783458956
